$wb = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item("meta")

# Copy the existing key-style (bold/orange) from A10 down to the new key cells
$meta.Range("A10").Copy() | Out-Null
$meta.Range("A11:A13").PasteSpecial(-4122) | Out-Null

$meta.Range("A11").Value = "x_date_format"
$meta.Range("B11").Value = "yyyy"

$meta.Range("A12").Value = "y_r_n_decimals"

# "1" looks numeric, so force text storage, then reset the cell format back
# to the plain/default style (matching the rest of column B) afterwards.
$meta.Range("B12").NumberFormat = "@"
$meta.Range("B12").Value = "1"
$meta.Range("B1").Copy() | Out-Null
$meta.Range("B12").PasteSpecial(-4122) | Out-Null

$meta.Range("A13").Value = $null
